# Adds the "ODI Bowling Extra" worksheet (with scraped bowling data) and
# cleans up the blank placeholder cells left behind in "ODI Batting Extra".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Tidy up "ODI Batting Extra": drop the leftover blank cells so only
#    cells that actually hold scraped data remain.
# ---------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$blankRanges = @(
    "B4:E4",
    "C5:E5",
    "B6:E6",
    "B7:E7",
    "B8:E8",
    "B9:E9",
    "C10:E10",
    "E11",
    "B12:E12",
    "C13:E13",
    "B14:E14",
    "B15:F21"
)
foreach ($r in $blankRanges) {
    $battingExtra.Range($r).ClearContents()
}

# ---------------------------------------------------------------------
# 2. Add the new "ODI Bowling Extra" sheet after "ODI Batting Extra".
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$bowlingExtra = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row, styled like the other "Extra" sheet's header.
$bowlingExtra.Cells.Item(1,1).Value = "MATCH_CODE"
$bowlingExtra.Cells.Item(1,2).Value = "MAIDEN_OVERS"
$bowlingExtra.Cells.Item(1,3).Value = "PERCENT_WICKETS_OF_ALL"

$battingExtra.Range("A1:C1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)

# Data rows.
$rows = @(
    @("4295", $null,  $null),
    @("4427", "3",    "30.00%"),
    @("4428", "0",    "20.00%"),
    @("4448", "1",    "10.00%"),
    @("4466", "1",    "30.00%"),
    @("4467", "1",    "40.00%"),
    @("4468", "0",    "10.00%"),
    @("4475", "0",    "20.00%"),
    @("4478", $null,  $null),
    @("4492", "0",    "10.00%"),
    @("4496", "0",    "30.00%"),
    @("4519", $null,  $null),
    @("4520", "1",    "20.00%"),
    @("4522", $null,  $null),
    @("4605", "0",    "10.00%"),
    @("4608", "0",    $null),
    @("4614", "0",    "20.00%"),
    @("4693", $null,  $null),
    @("4694", "0",    "40.00%"),
    @("4696", $null,  $null)
)

$r = 2
foreach ($row in $rows) {
    $bowlingExtra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) { $bowlingExtra.Cells.Item($r, 2).Value = $row[1] }
    if ($row[2] -ne $null) { $bowlingExtra.Cells.Item($r, 3).Value = $row[2] }
    $r = $r + 1
}

# Leave selection/activation where it started instead of on the new sheet.
$wb.Worksheets.Item("Player Info").Activate()
